$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.070.27'
$ws.Range('E2').Value = '  -0.61%  '
$ws.Range('D3').Value = '2.468.89'
$ws.Range('E3').Value = '  -0.43%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '581.83'
$ws.Range('E5').Value = '  -1.23%  '
$ws.Range('D6').Value = '168.23'
$ws.Range('E6').Value = '  -2.68%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '0.513'
$ws.Range('E8').Value = '  -1.67%  '
$ws.Range('D9').Value = '2.468.18'
$ws.Range('E9').Value = '  -0.45%  '
$ws.Range('E10').Value = '  -2.93%  '
$ws.Range('D11').Value = '0.164'
$ws.Range('E11').Value = '  -0.68%  '
$ws.Range('D12').Value = '4.95'
$ws.Range('E12').Value = '  -2.64%  '
$ws.Range('D13').Value = '0.332'
$ws.Range('E13').Value = '  -1.96%  '
$ws.Range('D14').Value = '25.49'
$ws.Range('E14').Value = '  -2.64%  '
$ws.Range('E15').Value = '  -1.19%  '
$ws.Range('D16').Value = '67.005.88'
$ws.Range('E16').Value = '  -0.54%  '
$ws.Range('E17').Value = '  -3.79%  '
$ws.Range('D18').Value = '2.472.72'
$ws.Range('E18').Value = '  -1.34%  '
$ws.Range('D19').Value = '11.24'
$ws.Range('E19').Value = '  -3.94%  '
$ws.Range('D20').Value = '7.52'
$ws.Range('E20').Value = '  -5.63%  '
$ws.Range('D21').Value = '355.79'
$ws.Range('E21').Value = '  -2.81%  '
$ws.Range('D22').Value = '4.03'
$ws.Range('E22').Value = '  -2.30%  '
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').Value = '69.15'
$ws.Range('E24').Value = '  -2.96%  '
$ws.Range('D25').Value = '4.23'
$ws.Range('E25').Value = '  -6.75%  '
$ws.Range('E26').Value = '  -6.29%  '
$ws.Range('D27').Value = '9.11'
$ws.Range('E27').Value = '  -8.09%  '
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  -0.05%  '
$ws.Range('D29').Value = '2.597.16'
$ws.Range('E29').Value = '  -0.98%  '
$ws.Range('D30').Value = '0.0₃0905'
$ws.Range('E30').Value = '  -5.43%  '
$ws.Range('D31').Value = '510.40'
$ws.Range('E31').Value = '  -3.76%  '
$ws.Range('D32').Value = '7.75'
$ws.Range('E32').Value = '  -7.51%  '
$ws.Range('E33').Value = '  -4.26%  '
$ws.Range('E34').Value = '  -5.26%  '
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('E36').Value = '  -6.49%  '
$ws.Range('D37').Value = '159.88'
$ws.Range('E37').Value = '  +1.17%  '
$ws.Range('D38').Value = '18.61'
$ws.Range('E38').Value = '  -0.08%  '
$ws.Range('D39').Value = '18.39'
$ws.Range('E39').Value = '  -1.57%  '
$ws.Range('E40').Value = '  -5.21%  '
$ws.Range('E41').Value = '  -0.34%  '
$ws.Range('E42').Value = '  -5.98%  '
$ws.Range('B43').Value = 'PolygonEcosystemToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D43').Value = '0.326'
$ws.Range('E43').Value = '  -6.34%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D44').Value = '4.79'
$ws.Range('E44').Value = '  -5.99%  '
$ws.Range('E45').Value = '  -5.66%  '
$ws.Range('D46').Value = '38.76'
$ws.Range('D47').Value = '140.85'
$ws.Range('E47').Value = '  -2.59%  '
$ws.Range('E48').Value = '  -5.53%  '
$ws.Range('D49').Value = '0.514'
$ws.Range('E49').Value = '  -5.72%  '
$ws.Range('E50').Value = '  -5.56%  '
$ws.Range('E51').Value = '  -8.84%  '
